$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert new row 4 (shifts old rows 4-25 down to 5-26) ---
$ws.Rows.Item(4).Copy()
$ws.Rows.Item(4).Insert()
$ws.Cells.Item(4, 1).Borders.LineStyle = 1

$ws.Cells.Item(4, 1).Value = 25
$ws.Cells.Item(4, 2).Value = 44538
$ws.Cells.Item(4, 3).Value = "推特迎新掌门人，印度为何「盛产」硅谷 CEO ？"
$ws.Cells.Item(4, 4).Value = "16:05"
$ws.Cells.Item(4, 5).Value = "https://aphid.fireside.fm/d/1437767933/12647593-905b-40ef-8977-371837f74e89/495fb43c-6a24-4b2d-a148-9f8291116987.mp3"

# --- Insert new row 22 (shifts old rows 22-26 down to 23-27) ---
$ws.Rows.Item(22).Copy()
$ws.Rows.Item(22).Insert()
$ws.Cells.Item(22, 1).Borders.LineStyle = 1

$ws.Cells.Item(22, 1).Value = 24
$ws.Cells.Item(22, 2).Value = 44531
$ws.Cells.Item(22, 3).Value = "The Big Consequences Of Small Changes To Congressional Maps"
$ws.Cells.Item(22, 4).Value = "894"
$ws.Cells.Item(22, 5).Value = "https://play.podtrac.com/npr-510310/edge1.pod.npr.org/anon.npr-mp3/npr/nprpolitics/2021/12/20211201_nprpolitics_120121politicspodcast.mp3?awCollectionId=510310&awEpisodeId=1060610347&orgId=1&d=894&p=510310&story=1060610347&t=podcast&e=1060610347&size=14304698&ft=pod&f=510310"
